# Add the "Networks" table entries (server/client NetworkEntity task list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column K, and set widths for new columns L and M
$ws.Columns.Item(11).ColumnWidth = 34.42578125
$ws.Columns.Item(12).ColumnWidth = 23
$ws.Columns.Item(13).ColumnWidth = 15

# Create a temporary named cell style for the green "DONE" status cells in this
# table (fill only, no border) - this produces the extra cellXfs entry, then we
# remove the named style again so only the underlying cell format remains.
$doneStyle = $wb.Styles.Add("TempNetworksDone")
$doneStyle.Interior.Color = 5287936   # RGB(0,176,80) "DONE" green

# Row 2
$ws.Range("K2").Value = "Merge projects"
$ws.Range("M2").Value = "DONE"
$ws.Range("M2").Style = "TempNetworksDone"

# Row 3
$ws.Range("K3").Value = "Get server to generate maze"
$ws.Range("M3").Value = "DONE"
$ws.Range("M3").Style = "TempNetworksDone"

# Row 4
$ws.Range("K4").Value = "Promote Client to server"
$ws.Range("L4").Value = "When a server disconnects, make a client the server"
$ws.Range("M4").Value = "TODO"
$ws.Range("M4").Interior.Color = 255   # RGB(255,0,0) "TODO" red

# Row 5
$ws.Range("K5").Value = "Get server to send maze back"
$ws.Range("M5").Value = "TODO"
$ws.Range("M5").Interior.Color = 255

# Row 6
$ws.Range("K6").Value = "Get client to print out maze"
$ws.Range("M6").Value = "TODO"
$ws.Range("M6").Interior.Color = 255

# Row 7
$ws.Range("K7").Value = "Create server and client classes"
$ws.Range("L7").Value = "Tidy up the code a bit"
$ws.Range("M7").Value = "DONE"
$ws.Range("M7").Style = "TempNetworksDone"

# Remove the temporary named style definition; the cells keep their resolved
# direct formatting (fill only, no border) via the newly created cellXfs entry.
$wb.Styles.Item("TempNetworksDone").Delete()

# Update selection to M7 (matches the author's final cursor position)
$ws.Range("M7").Select()
